$wb = $excel.ActiveWorkbook

# Add the new "settings" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "settings"

# Populate data
$newSheet.Range("A1").Value = "tag label"
$newSheet.Range("B1").Value = "tag value"
$newSheet.Range("A3").Value = "label1"
$newSheet.Range("B3").Value = "value1"
$newSheet.Range("A4").Value = "label2"
$newSheet.Range("B4").Value = "value2"
